$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, [string]$value) {
    # Force the cell to remain a text value (avoids Excel auto-converting
    # numeric-looking strings like "243.67" or "11" into real numbers).
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

# Row 2
Set-TextValue $ws.Range("D2") "243.67"
Set-TextValue $ws.Range("G2") "11"

# Row 3
Set-TextValue $ws.Range("D3") "23.04"
Set-TextValue $ws.Range("G3") "11"

# Row 4
Set-TextValue $ws.Range("D4") "5.391"
Set-TextValue $ws.Range("G4") "11"

# Row 5
Set-TextValue $ws.Range("D5") "0.05947"
Set-TextValue $ws.Range("G5") "11"

# Row 6
Set-TextValue $ws.Range("G6") "11"

# Row 7
Set-TextValue $ws.Range("D7") "6.508"
Set-TextValue $ws.Range("G7") "11"

# Row 8
Set-TextValue $ws.Range("D8") "0.8105"
Set-TextValue $ws.Range("G8") "11"

# Row 9
Set-TextValue $ws.Range("D9") "0.9278"
Set-TextValue $ws.Range("G9") "11"

# Row 10
Set-TextValue $ws.Range("D10") "0.1430"
Set-TextValue $ws.Range("G10") "11"

# Row 11
Set-TextValue $ws.Range("D11") "0.07392"
Set-TextValue $ws.Range("G11") "11"

# Row 12
Set-TextValue $ws.Range("D12") "0.03268"
Set-TextValue $ws.Range("G12") "11"

# Row 13
Set-TextValue $ws.Range("D13") "0.03077"
Set-TextValue $ws.Range("G13") "11"

# Row 14
Set-TextValue $ws.Range("D14") "0.09360"
Set-TextValue $ws.Range("G14") "11"

# Row 15
Set-TextValue $ws.Range("G15") "11"

# Row 16
Set-TextValue $ws.Range("D16") "0.001570"
Set-TextValue $ws.Range("G16") "11"

# Row 17
Set-TextValue $ws.Range("D17") "0.04693"
Set-TextValue $ws.Range("G17") "11"

# Row 18
$ws.Range("B18").Value = "TigerCash"
$ws.Range("C18").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextValue $ws.Range("D18") "0.005948"
$ws.Range("E18").Value = "17TigerCashTCH"
Set-TextValue $ws.Range("G18") "11"

# Row 19
$ws.Range("B19").Value = "BitKan"
$ws.Range("C19").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
Set-TextValue $ws.Range("D19") "0.001256"
$ws.Range("E19").Value = "18BitKanKAN"
Set-TextValue $ws.Range("G19") "11"

# Row 20
$ws.Range("B20").Value = "HotbitToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
Set-TextValue $ws.Range("D20") "0.004786"
$ws.Range("E20").Value = "19HotbitTokenHTB"
Set-TextValue $ws.Range("G20") "11"

# Row 21
$ws.Range("B21").Value = "NitroEx"
$ws.Range("C21").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
Set-TextValue $ws.Range("D21") "0.00008001"
$ws.Range("E21").Value = "20NitroExNTX"
Set-TextValue $ws.Range("G21") "11"

# Row 22
$ws.Range("B22").Value = "LEO"
$ws.Range("C22").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue $ws.Range("D22") "3.566"
$ws.Range("E22").Value = "21LEOLEO"
Set-TextValue $ws.Range("G22") "11"

# Row 23
$ws.Range("B23").Value = "BTSEToken"
$ws.Range("C23").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
Set-TextValue $ws.Range("D23") "2.133"
$ws.Range("E23").Value = "22BTSETokenBTSE"
Set-TextValue $ws.Range("G23") "11"

# Row 24
$ws.Range("B24").Value = "One"
$ws.Range("C24").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextValue $ws.Range("D24") "0.01118"
$ws.Range("E24").Value = "23OneONEBestin24h"
Set-TextValue $ws.Range("G24") "11"

# Row 25
Set-TextValue $ws.Range("G25") "11"

# Row 26
Set-TextValue $ws.Range("D26") "0.1330"
Set-TextValue $ws.Range("G26") "11"

# Row 27
Set-TextValue $ws.Range("D27") "0.0002340"
Set-TextValue $ws.Range("G27") "11"

# Row 28
Set-TextValue $ws.Range("G28") "11"

# Row 29
Set-TextValue $ws.Range("G29") "11"

# Row 30
Set-TextValue $ws.Range("G30") "11"

# Row 31
Set-TextValue $ws.Range("G31") "11"

# Row 32
Set-TextValue $ws.Range("G32") "11"

# Row 33
Set-TextValue $ws.Range("G33") "11"

# Row 34
Set-TextValue $ws.Range("G34") "11"

# Row 35
Set-TextValue $ws.Range("G35") "11"

# Row 36
Set-TextValue $ws.Range("G36") "11"

# Row 37
Set-TextValue $ws.Range("G37") "11"

# Row 38
Set-TextValue $ws.Range("G38") "11"

# Row 39
Set-TextValue $ws.Range("G39") "11"

# Row 40
Set-TextValue $ws.Range("D40") "0.03926"
Set-TextValue $ws.Range("G40") "11"

# Row 41
Set-TextValue $ws.Range("D41") "0.006301"
Set-TextValue $ws.Range("G41") "11"

# Row 42
Set-TextValue $ws.Range("G42") "11"

# Row 43
Set-TextValue $ws.Range("D43") "0.003500"
Set-TextValue $ws.Range("G43") "11"

# Row 44
Set-TextValue $ws.Range("D44") "0.008832"
Set-TextValue $ws.Range("G44") "11"

# Row 45
Set-TextValue $ws.Range("D45") "0.00005166"
Set-TextValue $ws.Range("G45") "11"

# Row 46
Set-TextValue $ws.Range("G46") "11"

# Row 47
Set-TextValue $ws.Range("D47") "0.6780"
Set-TextValue $ws.Range("G47") "11"

# Row 48
Set-TextValue $ws.Range("G48") "11"

# Row 49
Set-TextValue $ws.Range("G49") "11"

# Row 50
Set-TextValue $ws.Range("G50") "11"

# Row 51
Set-TextValue $ws.Range("G51") "11"
